# Insert a new weekly price record for "Palta" (avocado) data.
# A new row is inserted before the current row 109, pushing the existing
# rows 109-129 down to 110-130 (dimension grows from T129 to T130).
# The newly inserted row 109 receives the new observation's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above row 109; everything below shifts down by one.
$ws.Rows.Item(109).Insert()

# Populate the newly inserted row 109 with the new record.
$ws.Cells.Item(109, 1).Value  = 1
$ws.Cells.Item(109, 2).Value  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(109, 3).Value  = "Arica y Parinacota"
$ws.Cells.Item(109, 4).Value  = 44816
$ws.Cells.Item(109, 5).Value  = 15
$ws.Cells.Item(109, 6).Value  = "Fruta"
$ws.Cells.Item(109, 7).Value  = 100106
$ws.Cells.Item(109, 8).Value  = "Oleaginosos"
$ws.Cells.Item(109, 9).Value  = 100106002
$ws.Cells.Item(109, 10).Value = "Palta"
$ws.Cells.Item(109, 11).Value = "Hass"
$ws.Cells.Item(109, 12).Value = "Primera"
$ws.Cells.Item(109, 13).Value = 400
$ws.Cells.Item(109, 14).Value = 21000
$ws.Cells.Item(109, 15).Value = 22000
$ws.Cells.Item(109, 16).Value = 21500
$ws.Cells.Item(109, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(109, 18).Value = "Perú"
$ws.Cells.Item(109, 19).Value = 2150
$ws.Cells.Item(109, 20).Value = 10
